{"js": "// The author's commit ran Word's grammar/proofing pass over the document\n// (resolving \"read pdf in server\"), which causes Word to:\n//   1) split a few runs around words it flags, wrapping the flagged word in\n//      <w:proofErr w:type=\"gramStart\"/> ... <w:proofErr w:type=\"gramEnd\"/>\n//      (no visible text change \u2014 same text, just re-run-ified), and\n//   2) add two literal tab characters (<w:tab/> runs + a matching tab stop)\n//      into a previously-empty bold paragraph.\n//\n// We reproduce each surgically via Range.insertOoxml(..., Replace) so only\n// the targeted run/paragraph is touched and all other formatting/numbering\n// on the paragraph is left intact.\n\nconst OOXML_NS =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">';\n\nfunction wrapOoxml(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    OOXML_NS +\n    '<w:body>' +\n    bodyInnerXml +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nasync function replaceRunTextWithGrammarSplit(context, searchText, before, flagged, after, rPrXml) {\n  // Finds the exact run text (searchText = before + flagged + after) and\n  // replaces it with 3 runs, wrapping the middle one in proofErr marks.\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  const r = results.items[0];\n  const runXml =\n    (before !== \"\"\n      ? `<w:r><w:rPr>${rPrXml}</w:rPr><w:t xml:space=\"preserve\">${before}</w:t></w:r>`\n      : \"\") +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    `<w:r><w:rPr>${rPrXml}</w:rPr><w:t>${flagged}</w:t></w:r>` +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    (after !== \"\"\n      ? `<w:r><w:rPr>${rPrXml}</w:rPr><w:t xml:space=\"preserve\">${after}</w:t></w:r>`\n      : \"\");\n  r.insertOoxml(wrapOoxml(`<w:p>${runXml}</w:p>`), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"Reviewed, edited the readme file and provided comments. \"\nawait replaceRunTextWithGrammarSplit(\n  context,\n  \"Reviewed, edited the readme file and provided comments. \",\n  \"Reviewed, edited the readme \",\n  \"file\",\n  \" and provided comments. \",\n  '<w:rFonts w:cs=\"Times New Roman (Body CS)\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/>'\n);\n\n// 2) \" at the developer level\"\nawait replaceRunTextWithGrammarSplit(\n  context,\n  \" at the developer level\",\n  \" at the developer \",\n  \"level\",\n  \"\",\n  '<w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/>'\n);\n\n// 3) \" as a user, and record the errors/problems/ imperfections.\"\nawait replaceRunTextWithGrammarSplit(\n  context,\n  \" as a user, and record the errors/problems/ imperfections.\",\n  \" as a \",\n  \"user, and\",\n  \" record the errors/problems/ imperfections.\",\n  '<w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/>'\n);\n\n// 4) \"Start working on a user manual\"\nawait replaceRunTextWithGrammarSplit(\n  context,\n  \"Start working on a user manual\",\n  \"Start working on a user \",\n  \"manual\",\n  \"\",\n  '<w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/>'\n);\n\n// 5) Insert two tab characters (as real <w:tab/> runs, with a matching tab\n// stop) into the previously-empty bold paragraph that sits right after\n// \"Develop a user manual that can be download by a user from the frontend.\"\n// and right before \"Specific tasks for the coming week\".\n{\n  const paras = context.document.body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n\n  let target = null;\n  for (let i = 0; i < paras.items.length; i++) {\n    const p = paras.items[i];\n    if (p.text === \"\" ) {\n      // confirm this is the right empty paragraph by checking neighbors\n      if (i > 0 && i + 1 < paras.items.length) {\n        const prevText = paras.items[i - 1].text;\n        const nextText = paras.items[i + 1].text;\n        if (\n          prevText.indexOf(\"Develop a user manual\") !== -1 &&\n          nextText.indexOf(\"Specific\") === 0\n        ) {\n          target = p;\n          break;\n        }\n      }\n    }\n  }\n  if (!target) {\n    throw new Error(\"Target empty paragraph not found\");\n  }\n\n  const pXml =\n    '<w:p w14:paraId=\"7202BB5A\" w14:textId=\"77777777\" w:rsidR=\"00A34745\" w:rsidRPr=\"006601F1\" w:rsidRDefault=\"00A34745\" w:rsidP=\"00A34745\" ' +\n    'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    '<w:pPr><w:tabs><w:tab w:val=\"left\" w:pos=\"6240\"/></w:tabs><w:ind w:left=\"720\" w:hanging=\"360\"/>' +\n    '<w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:tab/></w:r>' +\n    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:tab/></w:r>' +\n    '</w:p>';\n  target.insertOoxml(wrapOoxml(pXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The author's commit ran Word's grammar/proofing pass over the document\n# (resolving \"read pdf in server\"), which causes Word to:\n#   1) split a few runs around words it flags, wrapping the flagged word in\n#      <w:proofErr w:type=\"gramStart\"/> ... <w:proofErr w:type=\"gramEnd\"/>\n#      (no visible text change -- same text, just re-run-ified), and\n#   2) add two literal tab characters (<w:tab/> runs + a matching tab stop)\n#      into a previously-empty bold paragraph.\n#\n# We reproduce each surgically: locate the whole paragraph by its (trimmed)\n# text, then replace that paragraph's full Range via Range.InsertXML with an\n# OOXML fragment that keeps every existing run/property byte-identical and\n# only re-shapes the target run(s).\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($doc, $i) {\n    if ($i -lt 1 -or $i -gt $doc.Paragraphs.Count) { return $null }\n    return $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n}\n\nfunction Find-ParaIndexExact($doc, $matchText) {\n    $n = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $n; $i++) {\n        $t = Get-ParaText $doc $i\n        if ($t -eq $matchText) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Find-EmptyParaBetween($doc, $prevNeedle, $nextNeedle) {\n    $n = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $n; $i++) {\n        $t = Get-ParaText $doc $i\n        if ($t -eq \"\") {\n            $pi = $i - 1\n            $ni = $i + 1\n            $prev = Get-ParaText $doc $pi\n            $next = Get-ParaText $doc $ni\n            if ($prev -like $prevNeedle -and $next -like $nextNeedle) {\n                return $i\n            }\n        }\n    }\n    return -1\n}\n\nfunction Set-ParagraphXml($doc, $paraIndex, $innerBodyXml) {\n    $r = $doc.Paragraphs.Item($paraIndex).Range\n    $xml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $r.InsertXML($xml)\n}\n\n# 1) \"Reviewed, edited the readme file and provided comments. \"\n$idx1 = Find-ParaIndexExact $d \"Reviewed, edited the readme file and provided comments. \"\nif ($idx1 -eq -1) { throw \"paragraph 1 not found\" }\n$p1 = '<w:p w14:paraId=\"42F6E69A\" w14:textId=\"375DD4D4\" w:rsidR=\"001B6B09\" w:rsidRDefault=\"001B6B09\" w:rsidP=\"00D21B84\">' +\n      '<w:pPr><w:pStyle w:val=\"ListParagraph1\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr><w:rPr><w:rFonts w:cs=\"Times New Roman (Body CS)\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr>' +\n      '<w:r><w:rPr><w:rFonts w:cs=\"Times New Roman (Body CS)\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">Reviewed, edited the readme </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:rPr><w:rFonts w:cs=\"Times New Roman (Body CS)\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>file</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:rPr><w:rFonts w:cs=\"Times New Roman (Body CS)\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> and provided comments. </w:t></w:r>' +\n      '</w:p>'\nSet-ParagraphXml $d $idx1 $p1\n\n# 2) \"Develop, upgrade and maintenance a READ.ME file at the developer level\"\n$idx2 = Find-ParaIndexExact $d \"Develop, upgrade and maintenance a READ.ME file at the developer level\"\nif ($idx2 -eq -1) { throw \"paragraph 2 not found\" }\n$p2 = '<w:p w14:paraId=\"51786984\" w14:textId=\"1C38DD76\" w:rsidR=\"00D21B84\" w:rsidRDefault=\"00E27A1D\" w:rsidP=\"001B6DB7\">' +\n      '<w:pPr><w:ind w:left=\"720\" w:hanging=\"360\"/><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Develop, u</w:t></w:r>' +\n      '<w:r w:rsidR=\"00D506A3\"><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>pgrade and maintenance a READ.ME file</w:t></w:r>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\"> at the developer </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>level</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '</w:p>'\nSet-ParagraphXml $d $idx2 $p2\n\n# 3) \"Conduct testing of iMedbot repeatly as a user, and record the errors/problems/ imperfections.\"\n$idx3 = Find-ParaIndexExact $d \"Conduct testing of iMedbot repeatly as a user, and record the errors/problems/ imperfections.\"\nif ($idx3 -eq -1) { throw \"paragraph 3 not found\" }\n$p3 = '<w:p w14:paraId=\"00E228DE\" w14:textId=\"51C5C2F2\" w:rsidR=\"00D21B84\" w:rsidRDefault=\"00D21B84\" w:rsidP=\"00D21B84\">' +\n      '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\">Conduct testing of </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>iMedbot</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>repeatly</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\"> as a </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>user, and</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\"> record the errors/problems/ imperfections.</w:t></w:r>' +\n      '</w:p>'\nSet-ParagraphXml $d $idx3 $p3\n\n# 4) \"Start working on a user manual\"\n$idx4 = Find-ParaIndexExact $d \"Start working on a user manual\"\nif ($idx4 -eq -1) { throw \"paragraph 4 not found\" }\n$p4 = '<w:p w14:paraId=\"102BCFF8\" w14:textId=\"7CCD1A26\" w:rsidR=\"00E27A1D\" w:rsidRDefault=\"00E27A1D\" w:rsidP=\"00D21B84\">' +\n      '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t xml:space=\"preserve\">Start working on a user </w:t></w:r>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r><w:rPr><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>manual</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      '</w:p>'\nSet-ParagraphXml $d $idx4 $p4\n\n# 5) Insert two tab characters (as real <w:tab/> runs, with a matching tab\n# stop) into the previously-empty bold paragraph that sits right after\n# \"Develop a user manual that can be download by a user from the frontend.\"\n# and right before \"Specific tasks for the coming week\".\n$idx5 = Find-EmptyParaBetween $d \"*Develop a user manual*\" \"Specific*\"\nif ($idx5 -eq -1) { throw \"paragraph 5 (empty tabs paragraph) not found\" }\n$p5 = '<w:p w14:paraId=\"7202BB5A\" w14:textId=\"77777777\" w:rsidR=\"00A34745\" w:rsidRPr=\"006601F1\" w:rsidRDefault=\"00A34745\" w:rsidP=\"00A34745\">' +\n      '<w:pPr><w:tabs><w:tab w:val=\"left\" w:pos=\"6240\"/></w:tabs><w:ind w:left=\"720\" w:hanging=\"360\"/><w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr>' +\n      '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:tab/></w:r>' +\n      '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:tab/></w:r>' +\n      '</w:p>'\nSet-ParagraphXml $d $idx5 $p5\n"}
